# "correção nos dados e inicio da analise PNAD 2009"
#
# The sheet had two empty section-header rows ("situação do domicílio" and
# "grandes regiões e unidades da federação") interleaved with the data rows.
# The fix removes those two header rows entirely (so every data row moves up
# and lines up with the correct label), which also drops the two now-unused
# strings from the shared-string table and shrinks the used range from
# A1:I40 down to A1:I38. It also relabels the first data-column header
# (currently the placeholder "unnamed: 1_level_1") as "total".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the placeholder column header to "total".
$ws.Range("B2").Value = "total"

# Remove the "situação do domicílio" header row (row 5) - it has no data,
# it's just a section label; removing it shifts "urbana"/"rural" (and
# everything below) up by one row.
$ws.Rows(5).Delete()

# Remove the "grandes regiões e unidades da federação" header row. Before
# the delete above it was row 8; after that delete everything shifted up by
# one, so it is now row 7.
$ws.Rows(7).Delete()
